$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "27.112.66"
Set-TextValue "E2" "  -2.50%  "
Set-TextValue "D3" "1.865.64"
Set-TextValue "E3" "  -2.15%  "
Set-TextValue "D4" "1.001"
Set-TextValue "E4" "  +0.06%  "
Set-TextValue "D5" "306.22"
Set-TextValue "E5" "  -2.11%  "
Set-TextValue "E6" "  +0.01%  "
Set-TextValue "D7" "0.5138"
Set-TextValue "E7" "  -1.83%  "
Set-TextValue "D8" "0.3761"
Set-TextValue "E8" "  -0.62%  "
Set-TextValue "E9" "  -1.23%  "
Set-TextValue "D10" "0.8900"
Set-TextValue "E10" "  -1.85%  "
Set-TextValue "D11" "20.69"
Set-TextValue "E11" "  -2.98%  "
Set-TextValue "D12" "0.07566"
Set-TextValue "E12" "  -1.20%  "
Set-TextValue "D13" "1.851.52"
Set-TextValue "E13" "  -2.99%  "
Set-TextValue "E14" "  -2.72%  "
Set-TextValue "D15" "89.60"
Set-TextValue "E15" "  -2.69%  "
Set-TextValue "E16" "  -0.02%  "
Set-TextValue "D17" "0.000008461"
Set-TextValue "E17" "  -2.92%  "
Set-TextValue "E18" "  -3.36%  "
Set-TextValue "D19" "0.9999"
Set-TextValue "E19" "  -0.02%  "
Set-TextValue "D20" "27.143.09"
Set-TextValue "E20" "  -2.49%  "
Set-TextValue "D21" "5.017"
Set-TextValue "E21" "  -2.73%  "
Set-TextValue "D22" "2.090.77"
Set-TextValue "E22" "  -4.09%  "
Set-TextValue "E23" "  -3.49%  "
Set-TextValue "D24" "6.448"
Set-TextValue "E24" "  -2.89%  "
Set-TextValue "D25" "1.839"
Set-TextValue "E25" "  -1.84%  "
Set-TextValue "D26" "146.35"
Set-TextValue "E26" "  -4.78%  "
Set-TextValue "E27" "  -2.25%  "
Set-TextValue "D28" "2.089"
Set-TextValue "E28" "  -3.89%  "
Set-TextValue "D29" "112.83"
Set-TextValue "E29" "  -1.73%  "
Set-TextValue "D30" "4.657"
Set-TextValue "E30" "  -4.21%  "
Set-TextValue "D31" "4.654"
Set-TextValue "E31" "  -4.26%  "
Set-TextValue "D32" "0.09123"
Set-TextValue "E32" "  +0.76%  "
Set-TextValue "E33" "  -3.27%  "
Set-TextValue "D34" "3.072"
Set-TextValue "E34" "  -3.42%  "
Set-TextValue "E35" "  -6.17%  "
Set-TextValue "D36" "0.7254"
Set-TextValue "E36" "  -7.35%  "
Set-TextValue "D37" "0.02035"
Set-TextValue "E37" "  -3.07%  "
Set-TextValue "D38" "3.085"
Set-TextValue "E38" "  +0.40%  "
Set-TextValue "D39" "2.493"
Set-TextValue "E39" "  -4.53%  "
Set-TextValue "E40" "  -1.66%  "
Set-TextValue "D41" "0.5281"
Set-TextValue "E41" "  -5.91%  "
Set-TextValue "D42" "6.469"
Set-TextValue "E42" "  -3.84%  "
Set-TextValue "D43" "115.69"
Set-TextValue "E43" "  +0.25%  "
Set-TextValue "D44" "8.275"
Set-TextValue "E44" "  -3.55%  "
Set-TextValue "E45" "  -3.57%  "
Set-TextValue "D46" "0.9999"
Set-TextValue "E46" "  -0.03%  "
Set-TextValue "D47" "0.4620"
Set-TextValue "E47" "  -4.17%  "
Set-TextValue "D48" "9.951"
Set-TextValue "E48" "  -5.25%  "
Set-TextValue "D49" "1.564"
Set-TextValue "E49" "  -3.59%  "
Set-TextValue "D50" "36.51"
Set-TextValue "E50" "  -1.56%  "
Set-TextValue "D51" "63.47"
Set-TextValue "E51" "  -5.15%  "
